$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,("chicago-bears", "washington-commanders", -105, -115, 1, 38, 6, 2022)
    ,("pittsburgh-steelers", "tampa-bay-buccaneers", 360, -455, 9.5, 46, 6, 2022)
    ,("new-orleans-saints", "cincinnati-bengals", 125, -140, 2.5, 43, 6, 2022)
    ,("indianapolis-colts", "jacksonville-jaguars", -130, 110, 1.5, 41.5, 6, 2022)
    ,("cleveland-browns", "new-england-patriots", -145, 125, 2.5, 43, 6, 2022)
    ,("green-bay-packers", "new-york-jets", -340, 280, 7.5, 45, 6, 2022)
    ,("new-york-giants", "baltimore-ravens", 205, -245, 5.5, 46.5, 6, 2022)
    ,("miami-dolphins", "minnesota-vikings", 145, -165, 3, 45.5, 6, 2022)
    ,("atlanta-falcons", "san-francisco-49ers", 180, -210, 4.5, 44.5, 6, 2022)
    ,("los-angeles-rams", "carolina-panthers", -450, 370, 10, 42, 6, 2022)
    ,("seattle-seahawks", "arizona-cardianls", 125, -145, 2.5, 50.5, 6, 2022)
    ,("kansas-city-chiefs", "buffalo-bills", 125, -145, 2.5, 54, 6, 2022)
    ,("philadelphia-eagles", "dallas-cowboys", -265, 225, 6, 42, 6, 2022)
    ,("los-angeles-chargers", "denver-broncos", -235, 195, 5, 45.5, 6, 2022)
    ,("arizona-cardianls", "new-orleans-saints", -130, 110, 1.5, 44, 7, 2022)
    ,("dallas-cowboys", "detroit-lions", -330, 240, 7.5, 49, 7, 2022)
    ,("jacksonville-jaguars", "new-york-giants", -160, 140, 3, 42, 7, 2022)
    ,("tennessee-titans", "indianapolis-colts", -145, 125, 2.5, 42.5, 7, 2022)
    ,("cincinnati-bengals", "atlanta-falcons", -265, 195, 6, 47.5, 7, 2022)
    ,("baltimore-ravens", "cleveland-browns", -275, 205, 6.5, 45.5, 7, 2022)
    ,("carolina-panthers", "tampa-bay-buccaneers", 330, -505, 11, 40.5, 7, 2022)
    ,("washington-commanders", "green-bay-packers", 175, -215, 5.5, 41.5, 7, 2022)
    ,("las-vegas-raiders", "houston-texans", -300, 220, 7, 45.5, 7, 2022)
    ,("san-francisco-49ers", "kansas-city-chiefs", 125, -145, 2.5, 48.5, 7, 2022)
    ,("los-angeles-chargers", "seattle-seahawks", -250, 190, 6, 51, 7, 2022)
    ,("miami-dolphins", "pittsburgh-steelers", -315, 225, 7, 45, 7, 2022)
    ,("new-england-patriots", "chicago-bears", -370, 270, 9, 40, 7, 2022)
    ,("jacksonville-jaguars", "denver-broncos", -135, 115, 2.5, 40.5, 8, 2022)
    ,("dallas-cowboys", "chicago-bears", -475, 375, 10, 43, 8, 2022)
    ,("new-orleans-saints", "las-vegas-raiders", 100, -120, 1.5, 49.5, 8, 2022)
    ,("atlanta-falcons", "carolina-panthers", -210, 175, 4, 41, 8, 2022)
    ,("philadelphia-eagles", "pittsburgh-steelers", -500, 405, 10.5, 42.5, 8, 2022)
    ,("detroit-lions-", "miami-dolphins", 160, -185, 3.5, 51.5, 8, 2022)
    ,("minnesota-vikings", "arizona-cardianls", -190, 165, 3.5, 49, 8, 2022)
    ,("new-york-jets", "new-england-patriots", 125, -145, 2.5, 40, 8, 2022)
    ,("houston-texans", "tennessee-titans", -110, -110, 0, 39, 8, 2022)
    ,("seattle-seahawks", "new-york-giants", -165, 145, 3, 44.5, 8, 2022)
    ,("indianapolis-colts", "washington-commanders", -150, 130, 3, 39.5, 8, 2022)
    ,("los-angeles-rams", "san-francisco-49ers", -110, -110, 0, 42, 8, 2022)
    ,("buffalo-bills", "green-bay-packers", -525, 410, 10.5, 47, 8, 2022)
    ,("cleveland-browns", "cincinnati-bengals", 150, -170, 3, 45, 8, 2022)
    ,("houston-texans", "philadelphia-eagles", 550, -800, 13.5, 45, 9, 2022)
)

$startRow = 82
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}

$ws.Range("A124").Select()
